$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update report title and rolling-12-months label (October -> November) ---
$ws.Range("A1").Value = "Table 1.2.E. Net Generation by Energy Source:  Residential Sector, 2014-November 2016"

# --- Insert a new row for the November monthly figure (row 45), pushing the
#     existing Year to Date / Rolling 12 months blocks down by one row ---
$ws.Rows("45").Insert()

# Copy the formatting of the October row (now row 44) into the new row so the
# new cells keep the same styles (right-aligned number-format text / number).
$ws.Range("A44:B44").Copy()
$ws.Range("A45:B45").PasteSpecial(-4122)

$ws.Range("A45").Value = "November"
$ws.Range("B45").Value = 717

# --- Update the "Rolling 12 Months Ending in October" label (now at row 50) ---
$ws.Range("A50").Value = "Rolling 12 Months Ending in November"

# --- Update Year to Date figures (rows 47-49 after the insert) ---
$ws.Range("B47").Value = 4585
$ws.Range("B48").Value = 6527
$ws.Range("B49").Value = 9851

# --- Update Rolling 12 Months Ending figures (rows 51-52 after the insert) ---
$ws.Range("B51").Value = 6890
$ws.Range("B52").Value = 10323

$excel.CutCopyMode = $false
